# Chemical list.xlsx - update CLASS classification table
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# Rows 10-15: update existing cells (style is already correct, s=10/11)
# ---------------------------------------------------------------
$ws.Range("A10").Value = "Remark : Classification of Hazardous Substances"

$ws.Range("B11").Value = "Class 1"
$ws.Range("C11").Value = "วัตถุระเบิด"

$ws.Range("B12").Value = "Class 2A"
$ws.Range("C12").Value = "ก๊าซอัด ก๊าซเหลว ก๊าซละลายได้ภายใต้ความดัน"

$ws.Range("B13").Value = "Class 2B"
$ws.Range("C13").Value = "ก๊าซภายใต้ความดันในกระป๋องเสปร์ย"

$ws.Range("B14").Value = "Class 3A"
$ws.Range("C14").Value = "ของเหลวไวไฟ จุดวาบไฟไม่เกิน 60 ºC"

$ws.Range("B15").Value = "Class 3B"
$ws.Range("C15").Value = "ของเหลวไวไฟ จุดวาบไฟมากกว่า 60 ºC-93 ºC คุณสมบัติเข้ากับน้ำไม่ได้"

$ws.Range("H10").Value = "Class 6.1A"
$ws.Range("I10").Value = "สารติดไฟได้ ที่มีคุณสมบัติเป็นพิษ"

$ws.Range("H11").Value = "Class 6.1B"
$ws.Range("I11").Value = "สารไม่ติดไฟ ที่มีคุณสมบัติเป็นพิษ"

$ws.Range("H12").Value = "Class 6.2"
$ws.Range("I12").Value = "สารติดเชื้อ"

$ws.Range("H13").Value = "Class 7"
$ws.Range("I13").Value = "สารกัมมันตรังสี"

$ws.Range("H14").Value = "Class 8A"
$ws.Range("I14").Value = "สารติดไฟ ที่มีคุณสมบัติกัดกร่อน"

$ws.Range("H15").Value = "Class 8B"
$ws.Range("I15").Value = "สารไม่ติดไฟ ที่มีคุณสมบัติกัดกร่อน"

# ---------------------------------------------------------------
# Row 16: new row, built off a copy of row 15 (keeps its s="10" cell
# styles/spans), then trim down to only B/C/H/I like the target.
# ---------------------------------------------------------------
$ws.Rows("15:15").Copy()
$ws.Rows("16:16").Insert(-4121)
$ws.Range("A16").Clear()
$ws.Range("D16:G16").Clear()
$ws.Range("J16:AC16").Clear()

$ws.Range("B16").Value = "Class 4.1A"
$ws.Range("H16").Value = "Class 9"

$ws.Range("C16").Font.Name = "FreesiaUPC"
$ws.Range("C16").Font.Size = 16
$ws.Range("C16").Value = "ของแข็งไวไฟ ที่มีคุณสมบัติระเบิด"

$ws.Range("I16").Font.Name = "FreesiaUPC"
$ws.Range("I16").Font.Size = 16
$ws.Range("I16").Value = "วัตถุอันตรายประเภทอื่นๆ"

# ---------------------------------------------------------------
# Rows 17-21: each new row is built from a copy of the row above it
# (which already only has B/C/H/I populated), so the new row also
# only gets those 4 cells and picks up the correct spans automatically.
# ---------------------------------------------------------------
$ws.Rows("16:16").Copy()
$ws.Rows("17:17").Insert(-4121)
$ws.Range("B17").Value = "Class 4.1B"
$ws.Range("C17").Value = "ของแข็งไวไฟ ที่ไม่มีคุณสมบัติระเบิด"
$ws.Range("H17").Value = "Class 10"
$ws.Range("I17").Value = "ของเหลวติดไฟ"

$ws.Rows("17:17").Copy()
$ws.Rows("18:18").Insert(-4121)
$ws.Range("B18").Value = "Class 4.2"
$ws.Range("C18").Value = "สารที่มีความเสี่ยงต่อการลุกไหม้ได้เอง"
$ws.Range("H18").Value = "Class 11"
$ws.Range("I18").Value = "ของแข็งติดไฟได้"

$ws.Rows("18:18").Copy()
$ws.Rows("19:19").Insert(-4121)
$ws.Range("B19").Value = "Class 4.3"
$ws.Range("C19").Value = "สารให้ก๊าซไวไฟ เมื่อสัมผัสกับน้ำ"
$ws.Range("H19").Value = "Class 12"
$ws.Range("I19").Value = "ของเหลวไม่ติดไฟ"

$ws.Rows("19:19").Copy()
$ws.Rows("20:20").Insert(-4121)
$ws.Range("B20").Value = "Class 5.1"
$ws.Range("C20").Value = "สารออกซิไดซ์"
$ws.Range("H20").Value = "Class 13"
$ws.Range("I20").Value = "ของแข็งไม่ติดไฟ"

$ws.Rows("20:20").Copy()
$ws.Rows("21:21").Insert(-4121)
$ws.Range("H21:I21").Clear()
$ws.Range("B21").Value = "Class 5.2"
$ws.Range("C21").Value = "สารเปอร์ออกซิไดซ์"

# ---------------------------------------------------------------
# View/selection + print setup
# ---------------------------------------------------------------
$ws.Range("H10:I20").Select()

$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
